# Rename the embedded logo pictures in the headers/footers.
#
# The BTec logo (alt text "BTec_Logo-Orange") used in both headers was
# saved internally as "image2.jpg" and needs to become "image1.jpg".
# The Pearson Edexcel logo (alt text ending "PearsonLogo.png") used in
# both footers was saved internally as "image1.png" and needs to become
# "image2.png".

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Sections(1).Headers.Count; $i++) {
    $hdr = $d.Sections(1).Headers($i)
    if ($hdr.Exists) {
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $shp = $hdr.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}

for ($i = 1; $i -le $d.Sections(1).Footers.Count; $i++) {
    $ftr = $d.Sections(1).Footers($i)
    if ($ftr.Exists) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $shp = $ftr.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}
